# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamp
# columns for the 8e6706c6-cdb4-4fb6-b910-12cd42910fde row across the
# Overview, zh-cn and de-de sheets, reflecting a newly regenerated report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 8e6706c6-... row
$wsOverview.Range("G4").Value = "2016-10-17 14:19:28"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-10-17 14:19:05"
$wsZhCn.Range("K4").Value = "2016-10-17 14:20:14"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H4").Value = "2016-10-17 14:19:28"
$wsDeDe.Range("K4").Value = "2016-10-17 14:20:53"
